# Revert capacity chart to show kilowatts (instead of watts) on the
# y-axis: the "Solar" column (E) on Sheet1, and the matching chart
# series, were entered in watts; this converts the stored values to
# kilowatts (divide by 1000), adjusts the shared number format so the
# values keep one decimal of precision, and updates the chart's value
# axis title/number format to match.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- 1. Convert the "Solar" (column E) capacity figures from watts to
#        kilowatts on the worksheet that backs the chart. ---
$solarRows = @{
    13 = 3
    14 = 0
    15 = 9.4
    16 = 16
    17 = 7.6
    18 = 2.9
    19 = 19.5
    20 = 15.2
    21 = 63.5
    22 = 13.6
    23 = 39.61
    24 = 43.28
    25 = 79.65000000000001
    26 = 25.038
}
foreach ($row in $solarRows.Keys) {
    $ws.Range("E$row").Value = $solarRows[$row]
}

# --- 2. The shared number format used by every data cell (B2:G26)
#        shows one decimal place now that the numbers are fractional
#        kilowatts rather than whole watts. ---
$ws.Range("B2:G26").NumberFormat = "#,##0.0"

# --- 3. Fix up the chart itself: axis title text and the value axis
#        number format (no more ">=1000 -> K" abbreviation, since the
#        numbers are already in kilowatts). ---
$chartObj = $ws.ChartObjects(1)
$chart = $chartObj.Chart
$valueAxis = $chart.Axes(2)  # xlValue
$valueAxis.AxisTitle.Text = "Kilowatts (kW)"
$valueAxis.TickLabels.NumberFormat = "#,##0"

# --- 4. Keep the chart's cached series values in sync with the
#        worksheet so the plotted bars reflect the new kilowatt
#        figures. ---
$solarSeries = $chart.SeriesCollection(4)  # Solar
$solarSeries.Values = $ws.Range("E2:E26")
$solarSeries.Formula = "=SERIES(Sheet1!`$E`$1,Sheet1!`$A`$2:`$A`$26,Sheet1!`$E`$2:`$E`$26,4)"
